$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New Name/Value config rows appended after the existing list (rows 17-24),
# following the same two-column "Name" / "Value" pattern used by rows 2-16.
$ws.Range("A17").Value = "strAFSColumn"
$ws.Range("B17").Value = "Company Code"

$ws.Range("A18").Value = "strAFSColumns"
$ws.Range("B18").Value = "Company Code,Vendor,Name 1,Indus.,Industry Name,Pmt. Method,PayT Accounting"

$ws.Range("A19").Value = "strFMSColumn"
$ws.Range("B19").Value = "C.Code"

$ws.Range("A20").Value = "strFMSColumns"
$ws.Range("B20").Value = "C.Code,Vendor,Name,Industry key,Industry Name,Payment Terms,Payment Methods"

$ws.Range("A21").Value = "strNewColumnName"
$ws.Range("B21").Value = "Vendor_Code"

$ws.Range("A22").Value = "strNewColumnName1"
$ws.Range("B22").Value = "Industry_Key_Name"

$ws.Range("A23").Value = "strRange"
$ws.Range("B23").Value = "A1"

$ws.Range("A24").Value = "strSheetName"
$ws.Range("B24").Value = "Sheet1"

# Reset formatting on the newly typed block back to the plain/default style.
$ws.Range("A17:B24").ClearFormats()

# Drop the trailing block of now-unneeded blank filler rows at the bottom of
# the sheet (the sheet used to pad all the way down to row 995).
$ws.Range("A989:Z995").EntireRow.Delete()

# Leave the selection where the user ended up after the edit.
$ws.Range("A16").Select() | Out-Null
